# Expand each "Détail(s)" list item into four list items:
#   Nombre, Détail(s), Mail, Numéro de téléphone
# (same numbered-list paragraph formatting reused for all of them)

$d = $word.ActiveDocument

# Collect the indices of every paragraph whose entire text is "Détail(s)"
# first (collection is 1-based); walk them back-to-front so inserting new
# paragraphs never invalidates indices we still have to process.
$targets = New-Object System.Collections.ArrayList
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Détail(s)`r") {
        [void]$targets.Add($i)
    }
}

for ($k = $targets.Count - 1; $k -ge 0; $k--) {
    $idx = $targets[$k]
    $p = $d.Paragraphs.Item($idx)

    # Insert "Mail" then "Numéro de téléphone" after the "Détail(s)" paragraph.
    $p.Range.InsertParagraphAfter()
    $mail = $d.Paragraphs.Item($idx + 1)
    $mail.Range.Text = "Mail"

    $mail.Range.InsertParagraphAfter()
    $phone = $d.Paragraphs.Item($idx + 2)
    $phone.Range.Text = "Numéro de téléphone"

    # Insert "Nombre" before the "Détail(s)" paragraph.
    $p.Range.InsertParagraphBefore()
    $nombre = $d.Paragraphs.Item($idx)
    $nombre.Range.Text = "Nombre"
}
